$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "'Volume 32   Number  14"

$ws.Range("C9").Value = "'Report Covering the Week  3/31/2025  Through  4/6/2025"

$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Value = "'0"

$ws.Range("E14").NumberFormat = "General"
$ws.Range("E14").Value = "'***.*"

$ws.Range("N14").Value = -60

$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 3

$ws.Range("E15").Value = 200

$ws.Range("F15").Value = 5

$ws.Range("G15").Value = 3

$ws.Range("H15").Value = 66.666666666666

$ws.Range("I15").Value = 15

$ws.Range("J15").Value = 9

$ws.Range("K15").Value = 66.666666666666

$ws.Range("L15").Value = 36.363636363636

$ws.Range("M15").Value = 114.285714285714

$ws.Range("N15").Value = -37.5

$ws.Range("C16").Value = 5

$ws.Range("D16").Value = 3

$ws.Range("E16").Value = 66.666666666666

$ws.Range("F16").Value = 26

$ws.Range("G16").Value = 15

$ws.Range("H16").Value = 73.333333333333

$ws.Range("I16").Value = 81

$ws.Range("J16").Value = 86

$ws.Range("K16").Value = -5.813953488372

$ws.Range("L16").Value = 24.615384615384

$ws.Range("M16").Value = -32.5

$ws.Range("N16").Value = -87.727272727272

$ws.Range("C17").Value = 11

$ws.Range("D17").Value = 14

$ws.Range("E17").Value = -21.428571428571

$ws.Range("F17").Value = 54

$ws.Range("G17").Value = 63

$ws.Range("H17").Value = -14.285714285714

$ws.Range("I17").Value = 178

$ws.Range("J17").Value = 193

$ws.Range("K17").Value = -7.772020725388

$ws.Range("L17").Value = 5.325443786982

$ws.Range("M17").Value = 83.505154639175

$ws.Range("N17").Value = -36.879432624113

$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("C18").Value = 3

$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 3

$ws.Range("E18").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E18").Value = 0

$ws.Range("F18").Value = 9

$ws.Range("G18").Value = 10

$ws.Range("H18").Value = -10

$ws.Range("I18").Value = 41

$ws.Range("J18").Value = 52

$ws.Range("K18").Value = -21.153846153846

$ws.Range("L18").Value = -21.153846153846

$ws.Range("M18").Value = -60.576923076923

$ws.Range("N18").Value = -92.100192678227

$ws.Range("C19").Value = 18

$ws.Range("D19").Value = 11

$ws.Range("E19").Value = 63.636363636363

$ws.Range("F19").Value = 61

$ws.Range("G19").Value = 41

$ws.Range("H19").Value = 48.780487804878

$ws.Range("I19").Value = 161

$ws.Range("J19").Value = 153

$ws.Range("K19").Value = 5.228758169934

$ws.Range("L19").Value = -6.936416184971

$ws.Range("M19").Value = 14.18439716312

$ws.Range("N19").Value = -15.263157894736

$ws.Range("D20").Value = 7

$ws.Range("E20").Value = -71.428571428571

$ws.Range("G20").Value = 22

$ws.Range("H20").Value = -50

$ws.Range("I20").Value = 40

$ws.Range("J20").Value = 52

$ws.Range("K20").Value = -23.076923076923

$ws.Range("L20").Value = -45.205479452054

$ws.Range("M20").Value = -48.051948051948

$ws.Range("N20").Value = -92.844364937388

$ws.Range("C21").Value = 42

$ws.Range("D21").Value = 39

$ws.Range("E21").Value = 7.692307692307

$ws.Range("F21").Value = 166

$ws.Range("G21").Value = 158

$ws.Range("H21").Value = 5.06329113924

$ws.Range("I21").Value = 520

$ws.Range("J21").Value = 550

$ws.Range("K21").Value = -5.454545454545

$ws.Range("L21").Value = -5.109489051094

$ws.Range("M21").Value = -5.62613430127

$ws.Range("N21").Value = -76.827094474153

$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 1

$ws.Range("E22").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E22").Value = -100

$ws.Range("J22").Value = 4

$ws.Range("K22").Value = -50

$ws.Range("C24").Value = 14

$ws.Range("D24").Value = 23

$ws.Range("E24").Value = -39.130434782608

$ws.Range("F24").Value = 76

$ws.Range("G24").Value = 119

$ws.Range("H24").Value = -36.134453781512

$ws.Range("I24").Value = 264

$ws.Range("J24").Value = 343

$ws.Range("K24").Value = -23.032069970845

$ws.Range("L24").Value = -24.571428571428

$ws.Range("M24").Value = 5.179282868525

$ws.Range("C25").NumberFormat = "#,##0"
$ws.Range("C25").Value = 3

$ws.Range("D25").Value = 1

$ws.Range("E25").Value = 200

$ws.Range("F25").Value = 6

$ws.Range("G25").Value = 21

$ws.Range("H25").Value = -71.428571428571

$ws.Range("I25").Value = 39

$ws.Range("J25").Value = 56

$ws.Range("K25").Value = -30.357142857142

$ws.Range("L25").Value = -52.439024390243

$ws.Range("C26").Value = 18

$ws.Range("D26").Value = 7

$ws.Range("E26").Value = 157.142857142857

$ws.Range("F26").Value = 77

$ws.Range("G26").Value = 62

$ws.Range("H26").Value = 24.193548387096

$ws.Range("I26").Value = 218

$ws.Range("J26").Value = 229

$ws.Range("K26").Value = -4.803493449781

$ws.Range("L26").Value = 4.807692307692

$ws.Range("M26").Value = 2.347417840375

$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 3

$ws.Range("D27").Value = 1

$ws.Range("E27").Value = 200

$ws.Range("F27").Value = 5

$ws.Range("G27").Value = 5

$ws.Range("H27").Value = 0

$ws.Range("I27").Value = 16

$ws.Range("J27").Value = 12

$ws.Range("K27").Value = 33.333333333333

$ws.Range("L27").Value = 14.285714285714

$ws.Range("C28").Value = 1

$ws.Range("D28").Value = 1

$ws.Range("E28").Value = 0

$ws.Range("F28").Value = 6

$ws.Range("G28").Value = 5

$ws.Range("H28").Value = 20

$ws.Range("I28").Value = 24

$ws.Range("J28").Value = 19

$ws.Range("K28").Value = 26.315789473684

$ws.Range("L28").Value = 20

$ws.Range("C29").NumberFormat = "General"
$ws.Range("C29").Value = "'0"

$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Value = "'0"

$ws.Range("E29").NumberFormat = "General"
$ws.Range("E29").Value = "'***.*"

$ws.Range("N29").Value = -83.720930232558

$ws.Range("C30").NumberFormat = "General"
$ws.Range("C30").Value = "'0"

$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Value = "'0"

$ws.Range("E30").NumberFormat = "General"
$ws.Range("E30").Value = "'***.*"

$ws.Range("N30").Value = -85

$ws.Range("C31").NumberFormat = "#,##0"
$ws.Range("C31").Value = 1

$ws.Range("F31").NumberFormat = "#,##0"
$ws.Range("F31").Value = 1

$ws.Range("I31").NumberFormat = "#,##0"
$ws.Range("I31").Value = 1

$ws.Range("K31").Value = -50
